# Apply "Trade #30 closed" update across the workbook.
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.38
$summary.Range("B4").Value = -0.62
$summary.Range("B5").Value = -0.41
$summary.Range("B6").Value = 30
$summary.Range("B7").Value = 7
$summary.Range("B9").Value = 23.33

# --- Strategy Status sheet (MarketMaking row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.38
$status.Range("D4").Value = 30
$status.Range("E4").Value = -0.62
$status.Range("F4").Value = -0.62
$status.Range("G4").Value = 23.33

# --- All Trades sheet (Trade #30, row 31) ---
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G31").Value = 0.96
$allTrades.Range("H31").Value = "CLOSED"
$allTrades.Range("I31").Value = 74.5455
$allTrades.Range("J31").Value = 0.41
$allTrades.Range("K31").Value = 99.38
$allTrades.Range("P31").Value = "early_exit"
$allTrades.Range("Q31").Value = 5.02

# --- MarketMaking sheet (Trade #30, row 31) ---
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Range("G31").Value = 0.96
$marketMaking.Range("H31").Value = "CLOSED"
$marketMaking.Range("I31").Value = 74.5455
$marketMaking.Range("J31").Value = 0.41
$marketMaking.Range("K31").Value = 99.38
$marketMaking.Range("P31").Value = "early_exit"
$marketMaking.Range("Q31").Value = 5.02
